# "defined data from outlier"
#
# Two FLT8 samples (FLT8-3 and FLT8-4) previously had clearly-outlying
# absolute-depth measurements on the "depth" sheet. Those two measurements
# were re-derived/corrected, and the corrected rows are now also included
# on the "depth_good" sheet (where they had been omitted as outliers
# before).

$wb = $excel.ActiveWorkbook

# --- 1. Correct the raw values on the "depth" sheet -----------------------
$wsDepth = $wb.Worksheets.Item("depth")
$wsDepth.Cells.Item(7, 5).Value = 1.428914   # E7 (FLT8-3): 10.802132 -> 1.428914
$wsDepth.Cells.Item(8, 5).Value = 2.092668   # E8 (FLT8-4): 15.16492983 -> 2.092668

# --- 2. Add the two corrected rows to the "depth_good" sheet --------------
$wsGood = $wb.Worksheets.Item("depth_good")

# Insert two blank rows above the current row 6, pushing every existing
# row (FLT8-5 ... LYDIT5-9) down by two rows.
$wsGood.Rows("6:7").Insert()

# Row 6: FLT8-3 / 45deg / cutting / Flint / 1.428914
$wsGood.Cells.Item(6, 1).Value = "FLT8-3"
$wsGood.Cells.Item(6, 2).Value = "45°"
$wsGood.Cells.Item(6, 3).Value = "cutting"
$wsGood.Cells.Item(6, 4).Value = "Flint"
$wsGood.Cells.Item(6, 5).Value = 1.428914

# Row 7: FLT8-4 / 35deg / cutting / Flint / 2.092668
$wsGood.Cells.Item(7, 1).Value = "FLT8-4"
$wsGood.Cells.Item(7, 2).Value = "35°"
$wsGood.Cells.Item(7, 3).Value = "cutting"
$wsGood.Cells.Item(7, 4).Value = "Flint"
$wsGood.Cells.Item(7, 5).Value = 2.092668
